# Problem #2, item 4/5: merge the split "you need / to pick 20 socks..."
# sentence back into one run, then add a new "5. A. ..." paragraph (with
# an explanatory answer) that inherits the relocated _GoBack bookmark,
# followed by a trailing blank paragraph - matching the author's
# "added explanation of solution to problem #2" edit.

$d = $word.ActiveDocument

# Locate the target paragraph ("4. A. Each solution somewhat ...") by its
# (still-unique) leading text rather than a hard-coded index.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("4. A. Each solution somewhat meets the goals")) {
        $target = $cand
        break
    }
}

$finalText = "4. A. Each solution somewhat meets the goals the only that works for both questions is that you need to pick 20 socks to guarantee 1 or 3 pairs."

# The paragraph currently holds two runs ("...you need " / "to pick 20
# socks...") split by a hidden _GoBack bookmark. Re-assigning identical
# text is treated as a no-op by this host, so first stamp a throw-away
# marker (forcing a genuine text change -> the run split collapses and
# the now-orphaned bookmark is dropped), then write the real text.
$r = $d.Range($target.Range.Start, $target.Range.End - 1)
$r.Text = "TEMP_MARKER_TEXT"
$target2 = $d.Paragraphs.Item($target.Index)
$r2 = $d.Range($target2.Range.Start, $target2.Range.End - 1)
$r2.Text = $finalText

# Blank paragraph right after it.
$target3 = $d.Paragraphs.Item($target2.Index)
$target3.Range.InsertParagraphAfter()

# New paragraph for item 5, starting as a single run with its full text;
# it gets split into "5. " / "A. The only guaranteed ..." below once the
# bookmark is re-inserted between them.
$blank = $d.Paragraphs.Item($target3.Index + 1)
$blank.Range.InsertParagraphAfter()

$p5 = $d.Paragraphs.Item($blank.Index + 1)
$p5Text = "5. A. The only guaranteed solution is to pick out 20 socks to get the pairs you need to answer the questions."
$r5 = $d.Range($p5.Range.Start, $p5.Range.End - 1)
$r5.Text = $p5Text

# Re-create the _GoBack bookmark collapsed right after "5. " (position 3
# within the paragraph), reproducing its original split-run placement.
$p5b = $d.Paragraphs.Item($p5.Index)
$bmPos = $p5b.Range.Start + 3
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))

# Trailing blank paragraph after the new item 5 paragraph.
$p5c = $d.Paragraphs.Item($p5b.Index)
$p5c.Range.InsertParagraphAfter()

Write-Output "done"
